$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Ngày cấp" column (C) now holds real dates instead of free-text strings.
# Apply the Vietnamese-style short date format used by the rest of the
# workbook's new cells, then write the actual date values (01/08/2025).
$ws.Range("C1:C3").NumberFormat = "dd/mm/yyyy;@"
$ws.Range("C2").Value = 45870
$ws.Range("C3").Value = 45870

# Move the active selection to C13, matching where the author left off editing.
$ws.Range("C13").Select()
